$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.867.18"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.663.45"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.96"
$ws.Range("E5").Value = "  +7.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3646"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.30"
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3273"
$ws.Range("E9").Value = "  +1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07093"
$ws.Range("E11").Value = "  +2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.090"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.54"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.662.04"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.622"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001049"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06682"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.65"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.946"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.85"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.61"
$ws.Range("E23").Value = "  +5.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.818.11"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.467"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.442"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.11"
$ws.Range("E27").Value = "  +3.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.73"
$ws.Range("E28").Value = "  +0.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.843.17"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.91"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.170"
$ws.Range("E31").Value = "  +8.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.070"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.714"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08475"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.650"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.21"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06232"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.228"
$ws.Range("E38").Value = "  +2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.183"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02280"
$ws.Range("E40").Value = "  +2.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2086"
$ws.Range("E41").Value = "  +2.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.272"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5957"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.52"
$ws.Range("E45").Value = "  +7.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.853"
$ws.Range("E46").Value = "  +3.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5669"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.78"
$ws.Range("E48").Value = "  +3.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.962"
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06979"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.193"
$ws.Range("E51").Value = "  +4.48%  "
